$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-5 get cyclically shifted: row5 -> row2, row2 -> row3,
# row3 -> row4, row4 -> row5 (columns A,B,D,E,F,G,H,Q,R carry the data).
# Capture current (before) values first so we don't clobber source data
# while writing the new layout.

$cols = @("A","B","D","E","F","G","H","Q","R")

$before = @{}
foreach ($r in 2..5) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

# New row 2 = old row 5
# New row 3 = old row 2
# New row 4 = old row 3
# New row 5 = old row 4
$mapping = @{ 2 = 5; 3 = 2; 4 = 3; 5 = 4 }

foreach ($destRow in 2..5) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $before[$srcRow][$c]
    }
}

# The empty "Bestämningsmetod" (AF) inline-string placeholder cell also
# travels with its row: it was on row 4 before the shift, and ends up on
# row 5 afterwards. Force an explicit empty-text value (not a blank/no
# value) on AF5, then strip the quote-prefix style that creates so the
# cell is a plain, unstyled empty text cell - matching AF4's original
# representation - before clearing AF4 itself.
$ws.Range("AF5").Value = "'"
$ws.Range("AF5").Style = "Normal"
$ws.Range("AF4").ClearContents()
